# ----------------------------------------------------------------------------
# Applies the edit described by the diff:
#  1. Fix sheet-name typo "Evaluación fianciera" -> "Evaluación financiera"
#     (workbook sheet tab, chart series formulas, defined name Print_Area).
#  2. Add a new sheet "Justificación Costos" at the end, with its content.
#  3. Update the AUTOSOL financial model: unit cost C12 formula and the
#     initial machinery investment B18, which ripple through the rest of
#     the "Evaluación financiera" sheet.
#  4. Minor selection/view tweaks.
# ----------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---- 1. Rename sheet 2 & fix references that don't auto-follow the rename ----
$wsFin = $wb.Worksheets.Item(2)
$oldName = $wsFin.Name
$wsFin.Name = "Evaluación financiera"

# Defined name "Print_Area" still points at the old sheet title text.
foreach ($n in $wb.Names) {
    if ($n.Name -like "*Print_Area*") {
        $n.RefersTo = "='Evaluación financiera'!`$A`$1:`$O`$34"
    }
}

# Chart series formulas embed the sheet name as literal text and don't
# follow a Worksheets.Name rename automatically, so patch them explicitly.
foreach ($co in $wsFin.ChartObjects()) {
    $chart = $co.Chart
    for ($i = 1; $i -le $chart.SeriesCollection().Count; $i++) {
        $s = $chart.SeriesCollection().Item($i)
        $s.Formula = ($s.Formula -replace [regex]::Escape($oldName), "Evaluación financiera")
    }
}

# ---- 2. Financial-model input changes on "Evaluación financiera" ----
$wsFin.Range("C12").Formula = "=(10*1400000+2400000)*'Proyección ventas'!L3"
$wsFin.Range("B18").Value = -951883006.27999997

# ---- 3. View/selection tweaks ----
$wsFin.Range("B28").Select() | Out-Null

$wsVentas = $wb.Worksheets.Item(1)
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1

# ---- 4. New sheet "Justificación Costos" ----
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsCost = $wb.Worksheets.Add($null, $lastSheet)
$wsCost.Name = "Justificación Costos"

# Column widths
$wsCost.Columns.Item(1).ColumnWidth = 16.109375
$wsCost.Columns.Item(2).ColumnWidth = 12.88671875
$wsCost.Columns.Item(7).ColumnWidth = 22.21875
$wsCost.Columns.Item(13).ColumnWidth = 13.77734375
$wsCost.Columns.Item(14).ColumnWidth = 12.44140625

# -- Row 1: section titles --
$wsCost.Range("A1:H1").Merge() | Out-Null
$wsCost.Range("A1").Value = "Costo materia prima y empaque"

$wsCost.Range("K1:M1").Merge() | Out-Null
$wsCost.Range("K1").Value = "Bodega y Mantenimientos"

$wsCost.Range("N1").Formula = "=M2+M3"

# -- Row 2 --
$wsCost.Range("A2").Value = "Costos promedio para fabricación de un juguete"
$wsCost.Range("B2").Formula = "=SUM(B3:B6)"

$wsCost.Range("D2").Value = "Costo promedio de los juguetes"
$wsCost.Range("E2").Formula = "=(SUM('Proyección ventas'!H3:H5)/3)"

$wsCost.Range("G2").Value = "Porcentaje de materia prima y empaque respecto al costo de venta promedio del juguete"
$wsCost.Range("H2").Formula = "=(B2/E2)*100"

$wsCost.Range("K2:L2").Merge() | Out-Null
$wsCost.Range("K2").Value = "Bodega de 300 m cuadrados en zona industrial de bogotá"
$wsCost.Range("M2").Formula = "=4000000"

# -- Row 3 --
$wsCost.Range("A3").Value = "Plástico de inyección"
$wsCost.Range("B3").Value = 750

$wsCost.Range("K3:L3").Merge() | Out-Null
$wsCost.Range("K3").Value = "Mantenimiento mensual estimado"
$wsCost.Range("M3").Value = 1000000

# -- Row 4 --
$wsCost.Range("A4").Value = "Bolsa de empaque"
$wsCost.Range("B4").Value = 100

# -- Row 5 --
$wsCost.Range("A5").Value = "Caja de empaque"
$wsCost.Range("B5").Formula = "=100/25"

# -- Row 6 --
$wsCost.Range("A6").Value = "Pintura"
$wsCost.Range("B6").Value = 1000

# -- Row 9: section title --
$wsCost.Range("A9:H9").Merge() | Out-Null
$wsCost.Range("A9").Value = "Energía eléctrica"

# -- Row 10 --
$wsCost.Range("A10").Value = "Gasto promedio de energía por fabricación de juguete"
$wsCost.Range("B10").Formula = "=E10*1.21"

$wsCost.Range("D10").Value = "Costo kV/h para industria en Bogotá"
$wsCost.Range("E10").Value = 305.81

$wsCost.Range("G10").Value = "Porcentaje de gasto de energía eléctrica respecto al costo promedio de venta del juguete"
$wsCost.Range("H10").Formula = "=(B10/E2)*100"

# -- Row heights --
$wsCost.Rows.Item(2).RowHeight = 70.8
$wsCost.Rows.Item(10).RowHeight = 72

# -- Formatting: title rows (A1:H1 / A9:H9) -- bold 18pt, centered, boxed --
$titleRanges = @("A1:H1", "A9:H9")
foreach ($addr in $titleRanges) {
    $r = $wsCost.Range($addr)
    $r.Font.Bold = $true
    $r.Font.Size = 18
    $r.HorizontalAlignment = -4108 # xlCenter
    $r.Borders.Item(7).LineStyle = 1
    $r.Borders.Item(10).LineStyle = 1
    $r.Borders.Item(8).LineStyle = 1
    $r.Borders.Item(9).LineStyle = 1
}

# K1:M1 title (bold 18pt, boxed, default/general alignment)
$r = $wsCost.Range("K1:M1")
$r.Font.Bold = $true
$r.Font.Size = 18
$r.Borders.Item(7).LineStyle = 1
$r.Borders.Item(10).LineStyle = 1
$r.Borders.Item(8).LineStyle = 1
$r.Borders.Item(9).LineStyle = 1

# N1 value cell: accounting-ish number with no decimals, centered, boxed
$r = $wsCost.Range("N1")
$r.NumberFormat = "_(* #,##0_);_(* \(#,##0\);_(* ""-""??_);_(@_)"
$r.Font.Bold = $true
$r.Font.Size = 18
$r.HorizontalAlignment = -4108 # xlCenter
$r.Borders.Item(7).LineStyle = 1
$r.Borders.Item(10).LineStyle = 1
$r.Borders.Item(8).LineStyle = 1
$r.Borders.Item(9).LineStyle = 1

# -- A2/D2/G2 (italic labels, wrap, vcenter, box w/o top) --
foreach ($addr in @("A2", "D2", "G2")) {
    $r = $wsCost.Range($addr)
    $r.Font.Bold = $true
    $r.VerticalAlignment = -4108 # xlCenter -> vertical center
    $r.WrapText = $true
    $r.Borders.Item(7).LineStyle = 1
    $r.Borders.Item(10).LineStyle = 1
    $r.Borders.Item(9).LineStyle = 1
}

# -- B2/E2/H2 (numeric results, box w/o top) --
$r = $wsCost.Range("B2")
$r.NumberFormat = "$#,##0.00"
$r.HorizontalAlignment = -4152 # xlRight
$r.VerticalAlignment = -4108  # xlCenter
$r.Borders.Item(7).LineStyle = 1
$r.Borders.Item(10).LineStyle = 1
$r.Borders.Item(9).LineStyle = 1

$r = $wsCost.Range("E2")
$r.NumberFormat = "$#,##0.00"
$r.VerticalAlignment = -4108 # xlCenter
$r.Borders.Item(7).LineStyle = 1
$r.Borders.Item(10).LineStyle = 1
$r.Borders.Item(9).LineStyle = 1

$r = $wsCost.Range("H2")
$r.HorizontalAlignment = -4108 # xlCenter
$r.VerticalAlignment = -4108   # xlCenter
$r.Borders.Item(7).LineStyle = 1
$r.Borders.Item(10).LineStyle = 1
$r.Borders.Item(9).LineStyle = 1

# -- K2:L2 / K3:L3 (existing-style boxed, bold, centered, wrap) --
foreach ($addr in @("K2:L2", "K3:L3")) {
    $r = $wsCost.Range($addr)
    $r.Font.Bold = $true
    $r.HorizontalAlignment = -4108 # xlCenter
    $r.WrapText = $true
    $r.Borders.Item(7).LineStyle = 1
    $r.Borders.Item(10).LineStyle = 1
    $r.Borders.Item(8).LineStyle = 1
    $r.Borders.Item(9).LineStyle = 1
}

# -- M2/M3 numeric, boxed --
$r = $wsCost.Range("M2")
$r.NumberFormat = "$#,##0.00"
$r.HorizontalAlignment = -4108 # xlCenter
$r.VerticalAlignment = -4108   # xlCenter
$r.Borders.Item(7).LineStyle = 1
$r.Borders.Item(10).LineStyle = 1
$r.Borders.Item(8).LineStyle = 1
$r.Borders.Item(9).LineStyle = 1

$r = $wsCost.Range("M3")
$r.NumberFormat = "_(* #,##0_);_(* \(#,##0\);_(* ""-""??_);_(@_)"
$r.HorizontalAlignment = -4108 # xlCenter
$r.Borders.Item(7).LineStyle = 1
$r.Borders.Item(10).LineStyle = 1
$r.Borders.Item(8).LineStyle = 1
$r.Borders.Item(9).LineStyle = 1

# -- A3 (italic, wrap, box) --
$r = $wsCost.Range("A3")
$r.Font.Italic = $true
$r.WrapText = $true
$r.Borders.Item(7).LineStyle = 1
$r.Borders.Item(10).LineStyle = 1
$r.Borders.Item(8).LineStyle = 1
$r.Borders.Item(9).LineStyle = 1

# -- A4/A5/A6 (italic, box) --
foreach ($addr in @("A4", "A5", "A6")) {
    $r = $wsCost.Range($addr)
    $r.Font.Italic = $true
    $r.Borders.Item(7).LineStyle = 1
    $r.Borders.Item(10).LineStyle = 1
    $r.Borders.Item(8).LineStyle = 1
    $r.Borders.Item(9).LineStyle = 1
}

# -- B3/B4/B5/B6 (numeric currency, right, vcenter, box) --
foreach ($addr in @("B3", "B4", "B5", "B6")) {
    $r = $wsCost.Range($addr)
    $r.NumberFormat = "$#,##0.00"
    $r.HorizontalAlignment = -4152 # xlRight
    $r.VerticalAlignment = -4108   # xlCenter
    $r.Borders.Item(7).LineStyle = 1
    $r.Borders.Item(10).LineStyle = 1
    $r.Borders.Item(8).LineStyle = 1
    $r.Borders.Item(9).LineStyle = 1
}

# -- A10 (bold, left, vcenter, wrap, box) --
$r = $wsCost.Range("A10")
$r.Font.Bold = $true
$r.HorizontalAlignment = -4131 # xlLeft
$r.VerticalAlignment = -4108   # xlCenter
$r.WrapText = $true
$r.Borders.Item(7).LineStyle = 1
$r.Borders.Item(10).LineStyle = 1
$r.Borders.Item(8).LineStyle = 1
$r.Borders.Item(9).LineStyle = 1

# -- B10/E10/H10 (numeric, centered, vcenter, box) --
foreach ($addr in @("B10", "E10", "H10")) {
    $r = $wsCost.Range($addr)
    $r.HorizontalAlignment = -4108 # xlCenter
    $r.VerticalAlignment = -4108   # xlCenter
    $r.Borders.Item(7).LineStyle = 1
    $r.Borders.Item(10).LineStyle = 1
    $r.Borders.Item(8).LineStyle = 1
    $r.Borders.Item(9).LineStyle = 1
}

# -- D10 (bold, center, vcenter, wrap, box) --
$r = $wsCost.Range("D10")
$r.Font.Bold = $true
$r.HorizontalAlignment = -4108 # xlCenter
$r.VerticalAlignment = -4108   # xlCenter
$r.WrapText = $true
$r.Borders.Item(7).LineStyle = 1
$r.Borders.Item(10).LineStyle = 1
$r.Borders.Item(8).LineStyle = 1
$r.Borders.Item(9).LineStyle = 1

$wsCost.Range("D20").Select() | Out-Null

# ---- 5. Recalculate everything & leave the financial sheet active ----
$excel.CalculateFull()
$wsFin.Activate()
